$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previous used range first (rows 1-8, cols A-B) so stale cells don't linger
$ws.Range("A1:B8").Clear()

# Write the new data set
$values = @(
    @("admin", "jaP#uv+QAp9l"),
    @("rt",    "jaP#uv+QAp9l"),
    @("abc",   "jaP#uv+QAp9l"),
    @("admin", "rtg"),
    @("admin", "abc"),
    @("dr",    "drt"),
    @("test",  "abc")
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 2).Value = $values[$i][1]
    $ws.Cells.Item($r, 1).Value = $values[$i][0]
}

# Remove the header highlight style so A1/B1 go back to default formatting
$ws.Range("A1:B1").Style = "Normal"

# Update the active selection to match the saved view
[void]$ws.Range("B4").Select()
